$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header for 07_02_2024
$ws.Range("G1").Value = "07_02_2024"

# Add new column values
$ws.Range("G2").Value = 899
$ws.Range("G3").Value = 807
$ws.Range("G4").Value = 1367
$ws.Range("G5").Value = 2866

# Select the last entered cell to mirror the saved selection state
$ws.Range("G5").Select()
